# Add a new quarter "2022-Q3" sheet right after "总计" and a matching
# summary row in "总计" (sheet1), shifting the existing quarters down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update the "总计" (summary) sheet: insert a new row for 2022-Q3 at
#    row 2, shifting the existing rows (2022-Q2 .. 2020-Q4) down by one.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# Cache the existing rows 2-7 (B/C/D) before they get overwritten.
$oldB = @{}
$oldC = @{}
$oldD = @{}
for ($r = 2; $r -le 7; $r++) {
    $oldB[$r] = $ws1.Cells.Item($r, 2).Value2
    $oldC[$r] = $ws1.Cells.Item($r, 3).Value2
    $oldD[$r] = $ws1.Cells.Item($r, 4).Value2
}

# Row 8 is brand new - grab the formatting (bold/centered/bordered style
# used by column A) from row 7 before we touch anything, by copying the
# whole cell (value + format) then overwriting the value afterwards.
$ws1.Cells.Item(7, 1).Copy($ws1.Cells.Item(8, 1))

# Shift rows 2-7 down into rows 3-8 (from bottom up).
for ($r = 7; $r -ge 2; $r--) {
    $nr = $r + 1
    $ws1.Cells.Item($nr, 1).Value = $r - 1
    $ws1.Cells.Item($nr, 2).Value = $oldB[$r]
    $ws1.Cells.Item($nr, 3).Value = $oldC[$r]
    $ws1.Cells.Item($nr, 4).Value = $oldD[$r]
}

# Write the new 2022-Q3 summary row into row 2.
$ws1.Cells.Item(2, 1).Value = 0
$ws1.Cells.Item(2, 2).Value = "2022-Q3"
$ws1.Cells.Item(2, 3).Value = 6
$ws1.Cells.Item(2, 4).Value = 0.71

# ---------------------------------------------------------------------
# 2. Insert a brand-new "2022-Q3" worksheet right after "总计" (i.e.
#    before the existing "2022-Q2" sheet), holding the fund holdings
#    detail for the new quarter.
# ---------------------------------------------------------------------
$new = $wb.Worksheets.Add($null, $ws1)
$new.Name = "2022-Q3"

# Borrow cell formatting (bold/centered header row + bordered column A)
# from an existing 7-row quarter sheet so the new sheet's styling
# matches its siblings exactly.
$template = $wb.Worksheets.Item("2021-Q2")
$template.Range("B1:H1").Copy($new.Range("B1:H1"))
$template.Range("A2:H7").Copy($new.Range("A2:H7"))

# Header row text.
$new.Range("B1").Value = "基金代码"
$new.Range("C1").Value = "基金名称"
$new.Range("D1").Value = "基金规模"
$new.Range("E1").Value = "股票总仓位"
$new.Range("F1").Value = "仓位占比"
$new.Range("G1").Value = "持有市值(亿元)"
$new.Range("H1").Value = "仓位排名"

# Fund holdings data: index, code, name, scale, total position, position
# ratio, held market value, position rank.
$data = @(
    @(0, "004702", "南方金融主题灵活配置混合A", "11.74", "92.15", "3.13", "0.3675", 10),
    @(1, "013500", "南方金融主题灵活配置混合C", "7.39", "92.15", "3.13", "0.2313", 10),
    @(2, "515760", "华夏中证浙江国资创新发展ETF", "2.04", "99.57", "4.97", "0.1014", 6),
    @(3, "090011", "大成核心双动力混合", "0.24", "92.56", "2.31", "0.0055", 10),
    @(4, "013590", "南方比较优势混合A", "0.46", "66.06", "0.75", "0.0034", 9),
    @(5, "013591", "南方比较优势混合C", "0.30", "66.06", "0.75", "0.0022", 9)
)

$r = 2
foreach ($row in $data) {
    $new.Cells.Item($r, 1).Value = $row[0]

    # Columns B, D, E, F, G hold numeric-looking text (fund codes, e.g.
    # "004702", keep their leading zero; percentages/amounts are stored
    # as plain text too, matching the rest of the workbook). Force text
    # with a leading apostrophe, then clear the auto-applied
    # "number stored as text" quote-prefix style so the cell format
    # matches the plain (un-styled) data cells used elsewhere.
    $new.Cells.Item($r, 2).Value = "'" + $row[1]
    $new.Cells.Item($r, 2).ClearFormats()

    $new.Cells.Item($r, 3).Value = $row[2]

    $new.Cells.Item($r, 4).Value = "'" + $row[3]
    $new.Cells.Item($r, 4).ClearFormats()

    $new.Cells.Item($r, 5).Value = "'" + $row[4]
    $new.Cells.Item($r, 5).ClearFormats()

    $new.Cells.Item($r, 6).Value = "'" + $row[5]
    $new.Cells.Item($r, 6).ClearFormats()

    $new.Cells.Item($r, 7).Value = "'" + $row[6]
    $new.Cells.Item($r, 7).ClearFormats()

    $new.Cells.Item($r, 8).Value = $row[7]
    $r++
}
